# chore(prd): remove prototype appendix, embedded images, and references
#
# Removes the "Appendix: Quick prototype" Heading 2 section (its heading,
# the "Figure: PDF page NN" captions, and the embedded screenshot images
# that follow it), leaving the "Appendix: Links" section untouched.

$d = $word.ActiveDocument

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $style = $p.Style.NameLocal
    $text = $p.Range.Text

    if ($startPara -eq $null -and $style -eq "Heading 2" -and $text.StartsWith("Appendix: Quick prototype")) {
        $startPara = $p
    }
    elseif ($startPara -ne $null -and $endPara -eq $null -and $style -eq "Heading 2" -and $text.StartsWith("Appendix: Links")) {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.Start)
    $delRange.Delete()
}
